# Update the 20x5 table of arithmetic problems in the document.
# The replacement values are listed in row-major (top-to-bottom, left-to-right) order,
# matching the order of cells in the table.

$values = @(
    "70-10=60",
    "32+52=84",
    "10+49=59",
    "58+11=69",
    "9-9=0",
    "55-43=12",
    "65-57=8",
    "20+0=20",
    "57-21=36",
    "56+16=72",
    "44+2=46",
    "55+31=86",
    "13+38=51",
    "44+47=91",
    "42+52=94",
    "1+86=87",
    "75-20=55",
    "8+67=75",
    "15+4=19",
    "99-11=88",
    "27-5=22",
    "16+32=48",
    "29-20=9",
    "65-21=44",
    "35-24=11",
    "87-31=56",
    "37-32=5",
    "14+15=29",
    "52-45=7",
    "59+16=75",
    "55-18=37",
    "73-21=52",
    "41+28=69",
    "33+44=77",
    "30+9=39",
    "42-41=1",
    "62+14=76",
    "56-24=32",
    "74-50=24",
    "46-16=30",
    "77-74=3",
    "29-10=19",
    "40+28=68",
    "9+15=24",
    "99-28=71",
    "5+67=72",
    "87-59=28",
    "34+25=59",
    "92-84=8",
    "22+24=46",
    "15+9=24",
    "48+10=58",
    "38+38=76",
    "44-17=27",
    "91-37=54",
    "57-10=47",
    "83-58=25",
    "57+3=60",
    "20+22=42",
    "73-28=45",
    "89-22=67",
    "57-0=57",
    "35-6=29",
    "29-22=7",
    "71-19=52",
    "19+48=67",
    "66+12=78",
    "31+0=31",
    "61-8=53",
    "81-54=27",
    "40+1=41",
    "11+80=91",
    "48-20=28",
    "81-38=43",
    "19+76=95",
    "38+16=54",
    "74+21=95",
    "1+69=70",
    "27+28=55",
    "23+67=90",
    "71-42=29",
    "91+0=91",
    "25-8=17",
    "39+23=62",
    "43+17=60",
    "90-80=10",
    "58+0=58",
    "20+63=83",
    "13+12=25",
    "40+42=82",
    "9-0=9",
    "54+36=90",
    "79+16=95",
    "46+35=81",
    "70-8=62",
    "99-81=18",
    "54-38=16",
    "15+48=63",
    "28-3=25",
    "50+34=84"
)

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$rows = $table.Rows.Count
$cols = $table.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $table.Cell($r, $c)
        $cellRange = $cell.Range
        # Trim the trailing cell-mark / paragraph-mark characters so only the
        # visible text is replaced, preserving run formatting.
        $cellRange.MoveEnd(12, -1) | Out-Null
        $cellRange.Text = $values[$idx]
        $idx++
    }
}
